$wb = $excel.ActiveWorkbook

# --- eval_cats_evaluable: swap "evaluate one model" / "compare internally" (rows 3 & 4) ---
$ws = $wb.Worksheets.Item("eval_cats_evaluable")
$ws.Range("A3").Value = "compare internally"
$ws.Range("A4").Value = "evaluate one model"

# --- eval_cats_evaluated: swap "evaluate one model" / "compare internally" (rows 3 & 4) ---
$ws = $wb.Worksheets.Item("eval_cats_evaluated")
$ws.Range("A3").Value = "compare internally"
$ws.Range("A4").Value = "evaluate one model"

# --- eval_metrics: reorder labels for rows 6-14 (counts/pcts stay put) ---
$ws = $wb.Worksheets.Item("eval_metrics")
$ws.Range("A6").Value = "CCC"
$ws.Range("A7").Value = "probabilistic interval performance metric used in the COVID-19 Forecast Hub"
$ws.Range("A8").Value = "R2"
$ws.Range("A10").Value = "MSE"
$ws.Range("A11").Value = "ARE"
$ws.Range("A12").Value = "MAPE"
$ws.Range("A13").Value = "MedAE"
$ws.Range("A14").Value = "AE"

# --- uncertainty_subcats: rotate labels for rows 2-6 ---
$ws = $wb.Worksheets.Item("uncertainty_subcats")
$ws.Range("A2").Value = "95% CIs"
$ws.Range("A3").Value = "no"
$ws.Range("A4").Value = "sensitivity analysis"
$ws.Range("A5").Value = "multiple CIs"
$ws.Range("A6").Value = "80% CIs"

# --- data_cats: swap "health risk factors" / "hospital resources" (rows 9 & 11) ---
$ws = $wb.Worksheets.Item("data_cats")
$ws.Range("A9").Value = "hospital resources"
$ws.Range("A11").Value = "health risk factors"
